$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New financial figures for rows 11-27, columns D:H (old "0" placeholders -> real figures)
$data = @{
    11 = @(2952322, 4403780, 5063791, 10656705, 14721432)
    12 = @(-2493603, -3602424, -4262854, -9201908, -13158971)
    13 = @(458719, 801356, 800937, 1454797, 1562461)
    14 = @(-66619, -84211, -99841, -177950, -265630)
    15 = @(0, 0, 0, 0, 0)
    16 = @(69800, -3027, -727, -9847, -9)
    17 = @(461900, 714118, 700369, 1267000, 1296822)
    18 = @(-64934, -107725, -21655, -32998, -34926)
    19 = @(180500, 24752, 79464, 20373, 242640)
    20 = @(577466, 631145, 758178, 1254375, 1504536)
    21 = @(-100224, -125541, -146048, -182297, -224314)
    22 = @(477242, 505604, 612130, 1072078, 1280222)
    23 = @(0, 0, 0, 0, 0)
    24 = @(477242, 505604, 612130, 1072078, 1280222)
    25 = @(2447, 2593, 3142, 5508, 6577)
    26 = @(195000, 195000, 194822, 194650, 194650)
    27 = @(2447, 2593, 3139, 5498, 6565)
}

$cols = @("D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
